$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 47: change B47 from text "3" to a true number 3, keep everything else the same.
$ws.Range("B47").Value = 3

# Row 48: new row of annotation data for Ying Tang.
$ws.Range("A48").Value = "Ying Tang"
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "3"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "无"
$ws.Range("D48").Value = "ACK"
$ws.Range("E48").Value = "WRI"
$ws.Range("F48").Value = "92b80f86-ee70-4a78-8469-1a9c33b052ed"
$ws.Range("G48").Value = "7Y52YHDS2X7ae_annotated.xlsx"
$ws.Range("H48").Value = "We will include this description in the new version of the paper."
